# md2ppt_demo.pptx — slide 11 (ER-diagram "Data Model" slide)
# Rotate the four entity boxes (USER / ORDER / PRODUCT / LINE_ITEM) one
# position around the layout, move the relationship-label chips to their
# new spots, and re-route the connector geometry/flips to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# ---- Entity headers (colored title bars) -------------------------------
$s.Shapes.Item(2).TextFrame.TextRange.Text  = "LINE_ITEM"   # was USER
$s.Shapes.Item(4).TextFrame.TextRange.Text  = "PRODUCT"     # was ORDER
$s.Shapes.Item(6).TextFrame.TextRange.Text  = "USER"        # was PRODUCT
$s.Shapes.Item(8).TextFrame.TextRange.Text  = "ORDER"       # was LINE_ITEM

# ---- Entity field lists --------------------------------------------------
$s.Shapes.Item(3).TextFrame.TextRange.Text = "int id PK
int order_id FK
int product_id FK
int quantity"

$s.Shapes.Item(5).TextFrame.TextRange.Text = "int id PK
string name
decimal price
int stock"

$s.Shapes.Item(7).TextFrame.TextRange.Text = "int id PK
string name
string email
date created_at"

$s.Shapes.Item(9).TextFrame.TextRange.Text = "int id PK
int user_id FK
date order_date
decimal total"

# ---- Relationship-label chips: reposition ("places"/"contains"/"ordered in") --
$lbl1 = $s.Shapes.Item(10)     # "places"
$lbl1.Left = 320.62993
$lbl1.Top  = 263.3071

$lbl2 = $s.Shapes.Item(11)     # "contains"
$lbl2.Left = 100.1575
$lbl2.Top  = 263.3071

$lbl3 = $s.Shapes.Item(12)     # "ordered in"
$lbl3.Left = 210.3937
$lbl3.Top  = 164.8819

# ---- Connectors: re-route geometry / flips -------------------------------
$conn61 = $s.Shapes.Item(13)
$conn61.HorizontalFlip = -1
$conn61.Left   = 52.9134
$conn61.Top    = 174.72441
$conn61.Width  = 614.1733
$conn61.Height = 196.8504

$conn62 = $s.Shapes.Item(14)
$conn62.HorizontalFlip = -1
$conn62.VerticalFlip   = -1
$conn62.Left   = 52.9134
$conn62.Top    = 174.72441
$conn62.Width  = 173.2284
$conn62.Height = 196.8504

$conn63 = $s.Shapes.Item(15)
$conn63.HorizontalFlip = -1
$conn63.Left   = 52.9134
$conn63.Top    = 174.72441
$conn63.Width  = 393.7008
$conn63.Height = 0
